# Peru Liga 1 - base update (29-05-2024 22:54)
# Several match rows had been entered in the wrong order; this script
# restores the correct row <-> record association by swapping/rotating
# the data columns (B:AD) between the affected rows. Column A (the
# running index) is left untouched since it simply numbers the rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $rng1 = $ws.Range("B$r1`:AD$r1")
    $rng2 = $ws.Range("B$r2`:AD$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

# Simple pairwise swaps
Swap-Rows 156 157
Swap-Rows 175 177
Swap-Rows 180 182
Swap-Rows 187 188
Swap-Rows 252 253
Swap-Rows 294 295
Swap-Rows 312 313

# Rotations (3 rows whose data shifts cyclically)
# Group 1: new183 = old185 ; new184 = old183 ; new185 = old184
$r183 = $ws.Range("B183:AD183").Value2
$r184 = $ws.Range("B184:AD184").Value2
$r185 = $ws.Range("B185:AD185").Value2
$ws.Range("B183:AD183").Value2 = $r185
$ws.Range("B184:AD184").Value2 = $r183
$ws.Range("B185:AD185").Value2 = $r184

# Group 2: new338 = old339 ; new339 = old340 ; new340 = old338
$r338 = $ws.Range("B338:AD338").Value2
$r339 = $ws.Range("B339:AD339").Value2
$r340 = $ws.Range("B340:AD340").Value2
$ws.Range("B338:AD338").Value2 = $r339
$ws.Range("B339:AD339").Value2 = $r340
$ws.Range("B340:AD340").Value2 = $r338
